# no-op test
$d = $word.ActiveDocument
